$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capitalize column headings in the table header row (row 6)
$ws.Range("A6").Value = " Kérdések száma"
$ws.Range("B6").Value = "Embedding  generálásai idő átlaga"
$ws.Range("C6").Value = "Kontextus összeállitási idő átlaga"
$ws.Range("D6").Value = "LLM feldolgozási idő átlaga"
$ws.Range("E6").Value = "Teljes feldoldozási idő átlaga"
$ws.Range("F6").Value = "Szemantikus hasonlóság mérékének  (BERTScore F1) átlaga (0-1) között"

# Update header / title text (A5)
$ws.Range("A5").Value = "küszöbérték: 0.90; LLM modell: gemini-2.0-flash; top_k=50"

# Update sheet view: scroll position and selection
$ws.Application.ActiveWindow.ScrollRow = 9
$ws.Range("B17:G21").Select()
